$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("G2").Value = "2016-08-13 23:31:03"
$ws2.Range("H2").Value = "2016-08-13 23:30:53"
$ws2.Range("K2").Value = "2016-08-13 23:31:25"
$ws3.Range("K2").Value = "2016-08-13 23:31:35"
